# Apply updated crypto price/volume figures (and one Cosmos/Toncoin row swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.526.53'
$ws.Range("E2").Value = '  +4.85%  '
$ws.Range("D3").Value = '2.502.46'
$ws.Range("E3").Value = '  +3.26%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''323.69'
$ws.Range("E5").Value = '  +1.81%  '
$ws.Range("D6").Value = '''108.65'
$ws.Range("E6").Value = '  +5.91%  '
$ws.Range("D7").Value = '''0.528'
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.548'
$ws.Range("E9").Value = '  +3.68%  '
$ws.Range("D10").Value = '''38.48'
$ws.Range("E10").Value = '  +8.22%  '
$ws.Range("D11").Value = '''0.0814'
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = '''18.54'
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").Value = '''7.22'
$ws.Range("E14").Value = '  +1.89%  '
$ws.Range("D15").Value = '2.880.57'
$ws.Range("E15").Value = '  +2.72%  '
$ws.Range("D16").Value = '2.479.61'
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").Value = '''0.851'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").Value = '47.411.54'
$ws.Range("E18").Value = '  +4.75%  '
$ws.Range("D19").Value = '''12.87'
$ws.Range("E19").Value = '  +5.43%  '
$ws.Range("D20").Value = '''6.61'
$ws.Range("E20").Value = '  +4.30%  '
$ws.Range("D21").Value = '0.0₃0942'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Value = '''70.89'
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("E23").Value = '  +6.79%  '
$ws.Range("D24").Value = '''252.27'
$ws.Range("E24").Value = '  +3.24%  '
$ws.Range("E25").Value = '  +3.93%  '
$ws.Range("D26").Value = '''26.34'
$ws.Range("E26").Value = '  +2.74%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '''10.10'
$ws.Range("E28").Value = '  +5.25%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.21'
$ws.Range("E29").Value = '  -3.28%  '
$ws.Range("D30").Value = '''35.45'
$ws.Range("E30").Value = '  +7.60%  '
$ws.Range("D31").Value = '''0.137'
$ws.Range("E31").Value = '  +9.77%  '
$ws.Range("D32").Value = '''49.49'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").Value = '''19.88'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").Value = '''5.43'
$ws.Range("E34").Value = '  +4.15%  '
$ws.Range("D35").Value = '''0.0790'
$ws.Range("E35").Value = '  +3.57%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '''1.99'
$ws.Range("E37").Value = '  +6.65%  '
$ws.Range("D38").Value = '''4.65'
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("D39").Value = '''3.00'
$ws.Range("E39").Value = '  +4.23%  '
$ws.Range("E40").Value = '  +2.15%  '
$ws.Range("D41").Value = '''122.23'
$ws.Range("E41").Value = '  -2.75%  '
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("D43").Value = '''21.33'
$ws.Range("E43").Value = '  +3.61%  '
$ws.Range("D44").Value = '''0.0298'
$ws.Range("E44").Value = '  +2.95%  '
$ws.Range("D45").Value = '1.970.50'
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("D46").Value = '''3.02'
$ws.Range("E46").Value = '  +3.59%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("D49").Value = '''9.17'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").Value = '''5.29'
$ws.Range("E50").Value = '  +10.80%  '
$ws.Range("D51").Value = '''79.90'
$ws.Range("E51").Value = '  +3.79%  '
